$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 previously held "Tags" with a one-off "blank" style; it now holds the
# lowercase "tags" with default (no) formatting. Clear the formatting first,
# then set the new value.
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = "tags"

# Reset the active selection back to A1 (the default), which drops the
# explicit <selection> element that pointed at C6.
$ws.Range("A1").Select()
